$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 36.130436
$ws.Range("I5").Value = 30
$ws.Range("K5").Value = 30
$ws.Range("M5").Value = 85
$ws.Range("H19").Value = 3832
$ws.Range("I19").Value = 3832
$ws.Range("K19").Value = 3832
$ws.Range("M19").Value = -3657
$ws.Range("H28").Value = 908.5
$ws.Range("I28").Value = 1120.3846
$ws.Range("K28").Value = 1120.3846
$ws.Range("M28").Value = -635.3846000000001
$ws.Range("H32").Value = 2771.6428
$ws.Range("I32").Value = 2875.6
$ws.Range("J32").Value = 2713.889
$ws.Range("K32").Value = 2875.6
$ws.Range("L32").Value = 2713.889
$ws.Range("M32").Value = -2549.6
$ws.Range("N32").Value = -3365.889
$ws.Range("H40").Value = 3036.3635
$ws.Range("I40").Value = 3187.5
$ws.Range("K40").Value = 3187.5
$ws.Range("M40").Value = -3012.5
$ws.Range("H51").Value = 73429.3
$ws.Range("J51").Value = 94827.71000000001
$ws.Range("L51").Value = 94827.71000000001
$ws.Range("N51").Value = -95795.71000000001
$ws.Range("H92").Value = 58823990
$ws.Range("I92").Value = 424.41666
$ws.Range("K92").Value = 424.41666
$ws.Range("M92").Value = 823.58334
$ws.Range("H96").Value = 998.5
$ws.Range("J96").Value = 998.5
$ws.Range("L96").Value = 2995.5
$ws.Range("N96").Value = -5741.5
$ws.Range("H99").Value = 450.57144
$ws.Range("I99").Value = 495.5
$ws.Range("K99").Value = 1486.5
$ws.Range("M99").Value = 11.5
$ws.Range("H100").Value = 2343.5264
$ws.Range("J100").Value = 3117.375
$ws.Range("L100").Value = 3117.375
$ws.Range("N100").Value = -4199.375
$ws.Range("H103").Value = 996
$ws.Range("I103").Value = 447.5
$ws.Range("J103").Value = 1080.3846
$ws.Range("K103").Value = 1342.5
$ws.Range("L103").Value = 3241.1538
$ws.Range("M103").Value = -756.5
$ws.Range("N103").Value = -4413.1538
$ws.Range("H107").Value = 11719920
$ws.Range("I107").Value = 5435804
$ws.Range("J107").Value = 27779326
$ws.Range("K107").Value = 5435804
$ws.Range("L107").Value = 27779326
$ws.Range("M107").Value = -5433884
$ws.Range("N107").Value = -27783166
$ws.Range("H109").Value = 49000
$ws.Range("J109").Value = 49000
$ws.Range("L109").Value = 49000
$ws.Range("N109").Value = -51774
$ws.Range("H113").Value = 50933444
$ws.Range("I113").Value = 22224466
$ws.Range("K113").Value = 22224466
$ws.Range("M113").Value = -22221212
$ws.Range("H116").Value = 22734500
$ws.Range("I116").Value = 41671500
$ws.Range("K116").Value = 41671500
$ws.Range("M116").Value = -41668058
$ws.Range("H129").Value = 1416.9333
$ws.Range("I129").Value = 542.5
$ws.Range("K129").Value = 1627.5
$ws.Range("M129").Value = 3372.5
$ws.Range("H131").Value = 2716.9333
$ws.Range("I131").Value = 2068.7273
$ws.Range("K131").Value = 6206.1819
$ws.Range("M131").Value = -1166.1819
$ws.Range("H132").Value = 3804
$ws.Range("I132").Value = 3560.111
$ws.Range("K132").Value = 10680.333
$ws.Range("M132").Value = -8150.332999999999
$ws.Range("H137").Value = 2534.76
$ws.Range("J137").Value = 2133.5
$ws.Range("L137").Value = 6400.5
$ws.Range("N137").Value = -11500.5
$ws.Range("H138").Value = 5709.9116
$ws.Range("I138").Value = 1923.2106
$ws.Range("J138").Value = 10506.4
$ws.Range("K138").Value = 5769.6318
$ws.Range("L138").Value = 31519.2
$ws.Range("M138").Value = -629.6318000000001
$ws.Range("N138").Value = -41799.2
$ws.Range("H141").Value = 3218.2
$ws.Range("I141").Value = 3023
$ws.Range("K141").Value = 9069
$ws.Range("M141").Value = -3889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 136.33333
$ws.Range("I5").Value = 187
$ws.Range("K5").Value = 187
$ws.Range("M5").Value = -75
$ws.Range("H32").Value = 1958033.2
$ws.Range("I32").Value = 2236438.2
$ws.Range("K32").Value = 2236438.2
$ws.Range("M32").Value = -2236151.2
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("H74").Value = 18028.691
$ws.Range("I74").Value = 44576.93
$ws.Range("K74").Value = 44576.93
$ws.Range("M74").Value = -43702.93
$ws.Range("H77").Value = 18028.691
$ws.Range("I77").Value = 44576.93
$ws.Range("K77").Value = 222884.65
$ws.Range("M77").Value = -218516.65
$ws.Range("H88").Value = 2604.125
$ws.Range("I88").Value = 1697
$ws.Range("K88").Value = 1697
$ws.Range("M88").Value = -1291
$ws.Range("H91").Value = 2604.125
$ws.Range("I91").Value = 1697
$ws.Range("K91").Value = 1697
$ws.Range("M91").Value = -293
$ws.Range("H97").Value = 7591184
$ws.Range("I97").Value = 1209.875
$ws.Range("K97").Value = 1209.875
$ws.Range("M97").Value = -713.875
$ws.Range("H126").Value = 5240.143
$ws.Range("I126").Value = 5240.143
$ws.Range("K126").Value = 15720.429
$ws.Range("M126").Value = -13250.429
$ws.Range("H132").Value = 4404.5386
$ws.Range("I132").Value = 3253.392
$ws.Range("K132").Value = 9760.175999999999
$ws.Range("M132").Value = -7230.175999999999
$ws.Range("H133").Value = 94825.71000000001
$ws.Range("J133").Value = 94825.71000000001
$ws.Range("L133").Value = 94825.71000000001
$ws.Range("N133").Value = -99885.71000000001
foreach ($addr in @("N53", "M63", "N63", "M66", "N66")) { $ws.Range($addr).ClearContents() }

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 136.33333
$ws.Range("I4").Value = 187
$ws.Range("K4").Value = 187
$ws.Range("M4").Value = -72
$ws.Range("H20").Value = 7938637.5
$ws.Range("I20").Value = 11112960
$ws.Range("J20").Value = 2832.1667
$ws.Range("K20").Value = 11112960
$ws.Range("L20").Value = 2832.1667
$ws.Range("M20").Value = -11112713
$ws.Range("N20").Value = -3326.1667
$ws.Range("H86").Value = 27781072
$ws.Range("I86").Value = 11907417
$ws.Range("J86").Value = 41670520
$ws.Range("K86").Value = 11907417
$ws.Range("L86").Value = 41670520
$ws.Range("M86").Value = -11906294
$ws.Range("N86").Value = -41672766
$ws.Range("H89").Value = 27781072
$ws.Range("I89").Value = 11907417
$ws.Range("J89").Value = 41670520
$ws.Range("K89").Value = 59537085
$ws.Range("L89").Value = 208352600
$ws.Range("M89").Value = -59531469
$ws.Range("N89").Value = -208363832
$ws.Range("H94").Value = 4644
$ws.Range("I94").Value = 1066.6666
$ws.Range("J94").Value = 10010
$ws.Range("K94").Value = 1066.6666
$ws.Range("L94").Value = 10010
$ws.Range("M94").Value = -615.6666
$ws.Range("N94").Value = -10912
$ws.Range("H105").Value = 50625.78
$ws.Range("I105").Value = 72001.82000000001
$ws.Range("J105").Value = 3598.5
$ws.Range("K105").Value = 72001.82000000001
$ws.Range("L105").Value = 3598.5
$ws.Range("M105").Value = -70254.82000000001
$ws.Range("N105").Value = -7092.5
$ws.Range("H134").Value = 5421.9346
$ws.Range("I134").Value = 2100.577
$ws.Range("J134").Value = 9739.700000000001
$ws.Range("K134").Value = 6301.731000000001
$ws.Range("L134").Value = 29219.1
$ws.Range("M134").Value = -3766.731000000001
$ws.Range("N134").Value = -34289.10000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = 50
$ws.Range("H31").Value = 9982
$ws.Range("I31").Value = 3851.5454
$ws.Range("J31").Value = 13047.228
$ws.Range("K31").Value = 3851.5454
$ws.Range("L31").Value = 13047.228
$ws.Range("M31").Value = -3556.5454
$ws.Range("N31").Value = -13637.228
$ws.Range("H34").Value = 9982
$ws.Range("I34").Value = 3851.5454
$ws.Range("J34").Value = 13047.228
$ws.Range("K34").Value = 3851.5454
$ws.Range("L34").Value = 13047.228
$ws.Range("M34").Value = -3649.5454
$ws.Range("N34").Value = -13451.228
$ws.Range("H58").Value = 19239258
$ws.Range("I58").Value = 55557720
$ws.Range("J58").Value = 11837.647
$ws.Range("K58").Value = 55557720
$ws.Range("L58").Value = 11837.647
$ws.Range("M58").Value = -55557517
$ws.Range("N58").Value = -12243.647
$ws.Range("H105").Value = 5957371.5
$ws.Range("J105").Value = 11999
$ws.Range("L105").Value = 11999
$ws.Range("N105").Value = -15493
$ws.Range("H107").Value = 1528.7037
$ws.Range("I107").Value = 1465.5
$ws.Range("J107").Value = 1579.2667
$ws.Range("K107").Value = 1465.5
$ws.Range("L107").Value = 1579.2667
$ws.Range("M107").Value = 454.5
$ws.Range("N107").Value = -5419.2667
$ws.Range("H122").Value = 1191.7142
$ws.Range("J122").Value = 1420.8572
$ws.Range("L122").Value = 4262.571599999999
$ws.Range("N122").Value = -9162.571599999999
$ws.Range("H132").Value = 4926.756
$ws.Range("I132").Value = 1576.3182
$ws.Range("J132").Value = 8806.210999999999
$ws.Range("K132").Value = 4728.9546
$ws.Range("L132").Value = 26418.633
$ws.Range("M132").Value = -2198.9546
$ws.Range("N132").Value = -31478.633
$ws.Range("H134").Value = 6404.2188
$ws.Range("I134").Value = 1558.5883
$ws.Range("J134").Value = 11895.934
$ws.Range("K134").Value = 4675.7649
$ws.Range("L134").Value = 35687.802
$ws.Range("M134").Value = -2140.7649
$ws.Range("N134").Value = -40757.802
$ws.Range("H136").Value = 19239258
$ws.Range("I136").Value = 55557720
$ws.Range("J136").Value = 11837.647
$ws.Range("K136").Value = 166673160
$ws.Range("L136").Value = 35512.94100000001
$ws.Range("M136").Value = -166670610
$ws.Range("N136").Value = -40612.94100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4562.5
$ws.Range("I56").Value = 4562.5
$ws.Range("K56").Value = 4562.5
$ws.Range("M56").Value = -4032.5
$ws.Range("H61").Value = 534.25
$ws.Range("I61").Value = 104.28571
$ws.Range("J61").Value = 1136.2
$ws.Range("K61").Value = 312.85713
$ws.Range("L61").Value = 3408.6
$ws.Range("M61").Value = -97.85712999999998
$ws.Range("N61").Value = -3838.6
$ws.Range("H131").Value = 1729.1945
$ws.Range("I131").Value = 702.8182
$ws.Range("J131").Value = 2180.8
$ws.Range("K131").Value = 2108.4546
$ws.Range("L131").Value = 6542.400000000001
$ws.Range("M131").Value = 2931.5454
$ws.Range("N131").Value = -16622.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1818322.6
$ws.Range("I2").Value = 69.85714
$ws.Range("J2").Value = 5000265
$ws.Range("K2").Value = 69.85714
$ws.Range("L2").Value = 5000265
$ws.Range("M2").Value = 43.14286
$ws.Range("N2").Value = -5000491
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("H13").Value = 2896.7144
$ws.Range("I13").Value = 94.333336
$ws.Range("J13").Value = 4998.5
$ws.Range("K13").Value = 94.333336
$ws.Range("L13").Value = 4998.5
$ws.Range("M13").Value = 44.666664
$ws.Range("N13").Value = -5276.5
$ws.Range("H14").Value = 56000000
$ws.Range("I14").Value = 56000000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 56000000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -55999832
$ws.Range("H17").Value = 346.15384
$ws.Range("I17").Value = 300
$ws.Range("K17").Value = 300
$ws.Range("M17").Value = -132
$ws.Range("H20").Value = 5000
$ws.Range("J20").Value = 7000
$ws.Range("L20").Value = 7000
$ws.Range("N20").Value = -7490
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("H33").Value = 40000
$ws.Range("J33").Value = 40000
$ws.Range("L33").Value = 40000
$ws.Range("N33").Value = -40504
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("H36").Value = 6559.6665
$ws.Range("I36").Value = 6559.6665
$ws.Range("K36").Value = 6559.6665
$ws.Range("M36").Value = -6074.6665
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("H41").Value = 77599.75
$ws.Range("J41").Value = 103333
$ws.Range("L41").Value = 103333
$ws.Range("N41").Value = -104043
$ws.Range("H43").Value = 2079.25
$ws.Range("I43").Value = 2079.25
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 2079.25
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -1928.25
$ws.Range("H44").Value = 10000
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("H47").Value = 40000
$ws.Range("J47").Value = 40000
$ws.Range("L47").Value = 40000
$ws.Range("N47").Value = -41136
$ws.Range("H52").Value = 89999.5
$ws.Range("J52").Value = 89999.5
$ws.Range("L52").Value = 89999.5
$ws.Range("N52").Value = -90517.5
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("H58").Value = 73285
$ws.Range("J58").Value = 73285
$ws.Range("L58").Value = 73285
$ws.Range("N58").Value = -73839
$ws.Range("H59").Value = 20000
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -21166
$ws.Range("H70").Value = 10419.917
$ws.Range("I70").Value = 10004
$ws.Range("K70").Value = 10004
$ws.Range("M70").Value = -9734
$ws.Range("H73").Value = 10419.917
$ws.Range("I73").Value = 10004
$ws.Range("K73").Value = 10004
$ws.Range("M73").Value = -9068
$ws.Range("H80").Value = 146569.42
$ws.Range("I80").Value = 3649.5
$ws.Range("K80").Value = 3649.5
$ws.Range("M80").Value = -2651.5
$ws.Range("H83").Value = 146569.42
$ws.Range("I83").Value = 3649.5
$ws.Range("K83").Value = 18247.5
$ws.Range("M83").Value = -13255.5
$ws.Range("H107").Value = 2619.6
$ws.Range("J107").Value = 1249.5
$ws.Range("L107").Value = 1249.5
$ws.Range("N107").Value = -5089.5
$ws.Range("H113").Value = 5983.5127
$ws.Range("I113").Value = 2605.75
$ws.Range("K113").Value = 2605.75
$ws.Range("M113").Value = -435.75
$ws.Range("H122").Value = 7939014
$ws.Range("I122").Value = 14286766
$ws.Range("J122").Value = 4324.25
$ws.Range("K122").Value = 42860298
$ws.Range("L122").Value = 12972.75
$ws.Range("M122").Value = -42857848
$ws.Range("N122").Value = -17872.75
$ws.Range("H126").Value = 2967.611
$ws.Range("I126").Value = 3046.2307
$ws.Range("K126").Value = 9138.6921
$ws.Range("M126").Value = -6668.6921
$ws.Range("H132").Value = 6058.143
$ws.Range("J132").Value = 11419.223
$ws.Range("L132").Value = 34257.669
$ws.Range("N132").Value = -39317.669
$ws.Range("H140").Value = 63499.5
$ws.Range("J140").Value = 63499.5
$ws.Range("L140").Value = 63499.5
$ws.Range("N140").Value = -73859.5
$ws.Range("H141").Value = 38356.5
$ws.Range("J141").Value = 41693.855
$ws.Range("L141").Value = 41693.855
$ws.Range("N141").Value = -52053.855
foreach ($addr in @("M9", "N14", "M22", "M35", "N35", "M40", "N40", "N43", "N44", "N53")) { $ws.Range($addr).ClearContents() }

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6119.1055
$ws.Range("I7").Value = 3806.6667
$ws.Range("K7").Value = 3806.6667
$ws.Range("M7").Value = -3694.6667
$ws.Range("H14").Value = 13250
$ws.Range("I14").Value = 14333.333
$ws.Range("K14").Value = 14333.333
$ws.Range("M14").Value = -14161.333
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("H55").Value = 43478692
$ws.Range("I55").Value = 142857220
$ws.Range("J55").Value = 587.8125
$ws.Range("K55").Value = 142857220
$ws.Range("L55").Value = 587.8125
$ws.Range("M55").Value = -142857047
$ws.Range("N55").Value = -933.8125
$ws.Range("H61").Value = 4702.5454
$ws.Range("I61").Value = 1750.8462
$ws.Range("K61").Value = 1750.8462
$ws.Range("M61").Value = -1548.8462
$ws.Range("H82").Value = 671961.9
$ws.Range("I82").Value = 1006928.6
$ws.Range("J82").Value = 2028.4286
$ws.Range("K82").Value = 1006928.6
$ws.Range("L82").Value = 2028.4286
$ws.Range("M82").Value = -1006567.6
$ws.Range("N82").Value = -2750.4286
$ws.Range("H85").Value = 671961.9
$ws.Range("I85").Value = 1006928.6
$ws.Range("J85").Value = 2028.4286
$ws.Range("K85").Value = 1006928.6
$ws.Range("L85").Value = 2028.4286
$ws.Range("M85").Value = -1005680.6
$ws.Range("N85").Value = -4524.4286
$ws.Range("H87").Value = 58626
$ws.Range("J87").Value = 58626
$ws.Range("L87").Value = 58626
$ws.Range("N87").Value = -60872
$ws.Range("H90").Value = 58626
$ws.Range("J90").Value = 58626
$ws.Range("L90").Value = 175878
$ws.Range("N90").Value = -187110
$ws.Range("H93").Value = 7884.9287
$ws.Range("I93").Value = 7299.625
$ws.Range("J93").Value = 8665.333000000001
$ws.Range("K93").Value = 7299.625
$ws.Range("L93").Value = 8665.333000000001
$ws.Range("M93").Value = -6051.625
$ws.Range("N93").Value = -11161.333
$ws.Range("H107").Value = 2560
$ws.Range("I107").Value = 2560
$ws.Range("K107").Value = 2560
$ws.Range("M107").Value = -640
$ws.Range("H113").Value = 4702.5454
$ws.Range("I113").Value = 1750.8462
$ws.Range("K113").Value = 1750.8462
$ws.Range("M113").Value = 419.1538
$ws.Range("H122").Value = 4562.1665
$ws.Range("I122").Value = 2624.6843
$ws.Range("K122").Value = 7874.0529
$ws.Range("M122").Value = -5424.0529
$ws.Range("H126").Value = 6119.1055
$ws.Range("I126").Value = 3806.6667
$ws.Range("K126").Value = 11420.0001
$ws.Range("M126").Value = -8950.000100000001
$ws.Range("H132").Value = 11911565
$ws.Range("I132").Value = 27780012
$ws.Range("J132").Value = 10228.917
$ws.Range("K132").Value = 83340036
$ws.Range("L132").Value = 30686.751
$ws.Range("M132").Value = -83337506
$ws.Range("N132").Value = -35746.751
$ws.Range("H136").Value = 10899.538
$ws.Range("I136").Value = 1849.5
$ws.Range("J136").Value = 12545
$ws.Range("K136").Value = 5548.5
$ws.Range("L136").Value = 37635
$ws.Range("M136").Value = -2998.5
$ws.Range("N136").Value = -42735
foreach ($addr in @("M22", "N22", "M27", "N27")) { $ws.Range($addr).ClearContents() }

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 762500
$ws.Range("J41").Value = 762500
$ws.Range("L41").Value = 762500
$ws.Range("N41").Value = -763280
$ws.Range("H54").Value = 15416.667
$ws.Range("I54").Value = 15000
$ws.Range("K54").Value = 15000
$ws.Range("M54").Value = -14480
$ws.Range("H62").Value = 2196.75
$ws.Range("I62").Value = 2149
$ws.Range("J62").Value = 2244.5
$ws.Range("K62").Value = 2149
$ws.Range("L62").Value = 2244.5
$ws.Range("M62").Value = -1525
$ws.Range("N62").Value = -3492.5
$ws.Range("H65").Value = 2196.75
$ws.Range("I65").Value = 2149
$ws.Range("J65").Value = 2244.5
$ws.Range("K65").Value = 10745
$ws.Range("L65").Value = 11222.5
$ws.Range("M65").Value = -7625
$ws.Range("N65").Value = -17462.5
$ws.Range("H122").Value = 4303.415
$ws.Range("I122").Value = 3622.9707
$ws.Range("K122").Value = 10868.9121
$ws.Range("M122").Value = -8418.9121
$ws.Range("H126").Value = 4002.4614
$ws.Range("I126").Value = 2690.4285
$ws.Range("K126").Value = 8071.2855
$ws.Range("M126").Value = -5601.2855
$ws.Range("H132").Value = 12831037
$ws.Range("I132").Value = 14289584
$ws.Range("J132").Value = 68750
$ws.Range("K132").Value = 42868752
$ws.Range("L132").Value = 206250
$ws.Range("M132").Value = -42866222
$ws.Range("N132").Value = -211310
$ws.Range("H136").Value = 41673880
$ws.Range("I136").Value = 166668180
$ws.Range("J136").Value = 9110.223
$ws.Range("K136").Value = 500004540
$ws.Range("L136").Value = 27330.669
$ws.Range("M136").Value = -500001990
$ws.Range("N136").Value = -32430.669
$ws.Range("H140").Value = 57500
$ws.Range("J140").Value = 57500
$ws.Range("L140").Value = 57500
$ws.Range("N140").Value = -67860
$ws.Range("H141").Value = 95000
$ws.Range("J141").Value = 95000
$ws.Range("L141").Value = 95000
$ws.Range("N141").Value = -105360
